# Apply the "Add files via upload" edit to Fragenbaum.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 17: rename "B1. Stammdaten & Betreiber" -> "B1. Stammkunde" and clear ParentID (B17)
# (set this first so the new shared string is allocated before the E2 attribute string)
$ws.Range("C17").Value = "B1. Stammkunde"
$ws.Range("B17").Clear()

# Row 2, column E: attribute list for Stammdaten gets ":pflicht" added to
# Name/Firma and Steuernummer attributes.
$ws.Range("E2").Value = "Name/Firma:text:pflicht;Anschrift:text;Steuernummer:text:pflicht;AMA-Betriebsnummer:text;SVS-Versicherungsnummer:text;Bankverbindung(IBAN/BIC):text;Familienstand:dropdown(ledig,verh.,geschieden,verwitwet);Kinder:number;Vollmacht&DSGVO-Einwilligung:checkbox;"

# Rows 22-24: clear ParentID cells (B column)
$ws.Range("B22").Clear()
$ws.Range("B23").Clear()
$ws.Range("B24").Clear()

# Update the selected cell (cosmetic, matches the saved view state)
$ws.Range("E2").Select()
